# "added update budget process"
# - Rename the worksheet from "Sheet1" to "category"
# - Remove the sample category rows (Salary/Food/Entertainment/Transport),
#   keeping only the header row (category_id / category_name)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "category"

# Drop the 4 data rows beneath the header, leaving just the header row intact.
$ws.Rows("2:5").Delete()
